# Update LR-pair TPM-derived metrics on Sheet1 (Egf-Erbb4) to match the
# refreshed TPM inputs: ligand/receptor expression, specificity, and edge
# weight columns (E-T) for rows 2-11. Numeric literals are written in plain
# decimal (not scientific notation) since the PS parser here doesn't accept
# exponent literals like 1.23E-05.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.1463976666666667
$ws.Range("H2").Value = 0.439193
$ws.Range("I2").Value = 0.1157910139257259
$ws.Range("J2").Value = 0.115791013925726
$ws.Range("P2").Value = 0.9810128591839974
$ws.Range("Q2").Value = 0.001507749569
$ws.Range("R2").Value = 0.013569746121
$ws.Range("S2").Value = 0.1135924736390905
$ws.Range("T2").Value = 0.1135924736390905

$ws.Range("G3").Value = 0.1463976666666667
$ws.Range("H3").Value = 0.439193
$ws.Range("I3").Value = 0.1157910139257259
$ws.Range("J3").Value = 0.115791013925726
$ws.Range("Q3").Value = 0.00002918193488888889
$ws.Range("R3").Value = 0.000262637414
$ws.Range("S3").Value = 0.002198540286635469
$ws.Range("T3").Value = 0.00219854028663547

$ws.Range("I4").Value = 0.4041732358198567
$ws.Range("J4").Value = 0.4041732358198568
$ws.Range("P4").Value = 0.9810128591839974
$ws.Range("S4").Value = 0.3964991416772857
$ws.Range("T4").Value = 0.3964991416772857

$ws.Range("I5").Value = 0.4041732358198567
$ws.Range("J5").Value = 0.4041732358198568
$ws.Range("S5").Value = 0.007674094142571021
$ws.Range("T5").Value = 0.007674094142571022

$ws.Range("E6").Value = 2
$ws.Range("F6").Value = 0.6666666666666666
$ws.Range("G6").Value = 0.1328766666666667
$ws.Range("H6").Value = 0.39863
$ws.Range("I6").Value = 0.1050967840589721
$ws.Range("J6").Value = 0.1050967840589721
$ws.Range("P6").Value = 0.9810128591839974
$ws.Range("Q6").Value = 0.00136849679
$ws.Range("R6").Value = 0.01231647111
$ws.Range("S6").Value = 0.1031012966207354
$ws.Range("T6").Value = 0.1031012966207354

$ws.Range("E7").Value = 2
$ws.Range("F7").Value = 0.6666666666666666
$ws.Range("G7").Value = 0.1328766666666667
$ws.Range("H7").Value = 0.39863
$ws.Range("I7").Value = 0.1050967840589721
$ws.Range("J7").Value = 0.1050967840589721
$ws.Range("Q7").Value = 0.00002648674888888889
$ws.Range("R7").Value = 0.00023838074
$ws.Range("S7").Value = 0.001995487438236714
$ws.Range("T7").Value = 0.001995487438236714

$ws.Range("G8").Value = 0.180116
$ws.Range("H8").Value = 0.5403480000000001
$ws.Range("I8").Value = 0.1424600182442301
$ws.Range("J8").Value = 0.1424600182442301
$ws.Range("P8").Value = 0.9810128591839974
$ws.Range("Q8").Value = 0.001855014684
$ws.Range("R8").Value = 0.016695132156
$ws.Range("S8").Value = 0.1397551098171766
$ws.Range("T8").Value = 0.1397551098171766

$ws.Range("G9").Value = 0.180116
$ws.Range("H9").Value = 0.5403480000000001
$ws.Range("I9").Value = 0.1424600182442301
$ws.Range("J9").Value = 0.1424600182442301
$ws.Range("Q9").Value = 0.00003590312266666667
$ws.Range("R9").Value = 0.000323128104
$ws.Range("S9").Value = 0.002704908427053488
$ws.Range("T9").Value = 0.002704908427053488

$ws.Range("G10").Value = 0.2939293333333333
$ws.Range("H10").Value = 0.881788
$ws.Range("I10").Value = 0.2324789479512151
$ws.Range("J10").Value = 0.2324789479512152
$ws.Range("P10").Value = 0.9810128591839974
$ws.Range("Q10").Value = 0.003027178204
$ws.Range("R10").Value = 0.027244603836
$ws.Range("S10").Value = 0.2280648374297093
$ws.Range("T10").Value = 0.2280648374297093

$ws.Range("G11").Value = 0.2939293333333333
$ws.Range("H11").Value = 0.881788
$ws.Range("I11").Value = 0.2324789479512151
$ws.Range("J11").Value = 0.2324789479512152
$ws.Range("Q11").Value = 0.00005858991377777778
$ws.Range("R11").Value = 0.000527309224
$ws.Range("S11").Value = 0.004414110521505846
$ws.Range("T11").Value = 0.004414110521505847
